$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure target cells keep their original text representation (not auto-converted
# to numbers/percentages by Excel) by forcing Text number format first.
$cells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'B8', 'C8', 'D8', 'E8', 'B9', 'C9', 'D9', 'E9', 'B10', 'C10', 'D10', 'E10', 'B11', 'C11', 'D11', 'E11', 'B12', 'C12', 'D12', 'E12', 'B13', 'C13', 'D13', 'E13', 'B14', 'C14', 'D14', 'E14', 'B15', 'C15', 'D15', 'E15', 'B16', 'C16', 'D16', 'E16', 'B17', 'C17', 'D17', 'E17', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'E50', 'E51')
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '305.44'
$ws.Range('E2').Value = '0.32%'
$ws.Range('D3').Value = '36.56'
$ws.Range('E3').Value = '2.62%'
$ws.Range('D4').Value = '5.021'
$ws.Range('E4').Value = '-1.44%'
$ws.Range('D5').Value = '0.07850'
$ws.Range('E5').Value = '0.24%'
$ws.Range('D6').Value = '2.160'
$ws.Range('E6').Value = '-3.69%'
$ws.Range('D7').Value = '8.038'
$ws.Range('E7').Value = '-1.07%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = '0.9225'
$ws.Range('E8').Value = '-0.36%'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Value = '0.09946'
$ws.Range('E9').Value = '2.04%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '0.1871'
$ws.Range('E10').Value = '2.66%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '0.08663'
$ws.Range('E11').Value = '-0.54%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.03599'
$ws.Range('E12').Value = '5.23%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.09932'
$ws.Range('E13').Value = '-0.02%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '0.001491'
$ws.Range('E14').Value = '0.72%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = '0.005679'
$ws.Range('E15').Value = '-0.47%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = '3.465'
$ws.Range('E16').Value = '-0.55%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = '4.052'
$ws.Range('E17').Value = '1.03%'
$ws.Range('E18').Value = '8.64%'
$ws.Range('D19').Value = '0.3449'
$ws.Range('E19').Value = '-0.30%'
$ws.Range('D20').Value = '0.1346'
$ws.Range('E20').Value = '1.85%'
$ws.Range('D21').Value = '4.924'
$ws.Range('E21').Value = '8.37%'
$ws.Range('D22').Value = '0.2202'
$ws.Range('E22').Value = '-1.62%'
$ws.Range('D23').Value = '0.04603'
$ws.Range('E23').Value = '-1.66%'
$ws.Range('D24').Value = '0.005185'
$ws.Range('E24').Value = '14.20%'
$ws.Range('D25').Value = '0.001233'
$ws.Range('E25').Value = '-0.95%'
$ws.Range('D26').Value = '0.0001401'
$ws.Range('E26').Value = '7.78%'
$ws.Range('D27').Value = '0.0002719'
$ws.Range('E27').Value = '0.73%'
$ws.Range('D39').Value = '0.01811'
$ws.Range('E39').Value = '3.02%'
$ws.Range('D40').Value = '0.04743'
$ws.Range('E40').Value = '0.55%'
$ws.Range('D41').Value = '0.007866'
$ws.Range('E41').Value = '-0.51%'
$ws.Range('D42').Value = '0.1408'
$ws.Range('E42').Value = '-0.90%'
$ws.Range('D43').Value = '0.007595'
$ws.Range('E43').Value = '-5.75%'
$ws.Range('D44').Value = '0.002222'
$ws.Range('E44').Value = '-3.40%'
$ws.Range('D45').Value = '0.01046'
$ws.Range('E45').Value = '14.64%'
$ws.Range('D46').Value = '0.00006554'
$ws.Range('E46').Value = '5.54%'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').Value = '0.05%'
$ws.Range('D48').Value = '0.0005804'
$ws.Range('E48').Value = '0.06%'
$ws.Range('D49').Value = '34.52'
$ws.Range('E49').Value = '502.25%'
$ws.Range('E50').Value = '0.11%'
$ws.Range('E51').Value = '0.05%'
